$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Bitcoin
$ws.Range("D2").Value = "'66.528.39"
$ws.Range("E2").Value = "  +1.02%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "'3.188.02"
$ws.Range("E3").Value = "  -0.32%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  +0.28%  "

# Row 5 - BNB
$ws.Range("D5").Value = "'602.87"
$ws.Range("E5").Value = "  +1.12%  "

# Row 6 - Solana
$ws.Range("D6").Value = "'155.91"
$ws.Range("E6").Value = "  +3.25%  "

# Row 7 - USDC
$ws.Range("E7").Value = "  +0.08%  "

# Row 8 - LidoStakedEther
$ws.Range("D8").Value = "'3.187.93"
$ws.Range("E8").Value = "  -0.42%  "

# Row 9 - XRP
$ws.Range("D9").Value = "'0.550"
$ws.Range("E9").Value = "  +2.84%  "

# Row 10 - Dogecoin
$ws.Range("E10").Value = "  -1.42%  "

# Row 11 - Toncoin
$ws.Range("E11").Value = "  -4.37%  "

# Row 12 - Cardano
$ws.Range("E12").Value = "  +0.47%  "

# Row 13 - ShibaInu
$ws.Range("E13").Value = "  -1.84%  "

# Row 14 - Avalanche
$ws.Range("D14").Value = "'38.99"
$ws.Range("E14").Value = "  +2.19%  "

# Row 15 - WrappedliquidstakedEther2.0
$ws.Range("D15").Value = "'3.710.02"
$ws.Range("E15").Value = "  +0.04%  "

# Row 16 - WrappedBTC
$ws.Range("D16").Value = "'66.554.21"
$ws.Range("E16").Value = "  +1.72%  "

# Row 17 - Polkadot
$ws.Range("D17").Value = "'7.39"
$ws.Range("E17").Value = "  +1.37%  "

# Row 18 - WrappedEther
$ws.Range("D18").Value = "'3.184.17"
$ws.Range("E18").Value = "  +0.09%  "

# Row 20 - BitcoinCash
$ws.Range("D20").Value = "'513.08"
$ws.Range("E20").Value = "  +0.25%  "

# Row 21 - Chainlink
$ws.Range("D21").Value = "'15.46"
$ws.Range("E21").Value = "  -2.46%  "

# Row 22 - Polygon
$ws.Range("D22").Value = "'0.735"
$ws.Range("E22").Value = "  +0.59%  "

# Row 23 - Uniswap
$ws.Range("D23").Value = "'8.15"
$ws.Range("E23").Value = "  +2.63%  "

# Row 24 - InternetComputer(DFINITY)
$ws.Range("D24").Value = "'14.96"
$ws.Range("E24").Value = "  -2.05%  "

# Row 25 - Litecoin
$ws.Range("E25").Value = "  -0.67%  "

# Row 26 - Dai
$ws.Range("D26").Value = "'0.999"
$ws.Range("E26").Value = "  -0.15%  "

# Row 27 - PancakeSwap
$ws.Range("D27").Value = "'3.01"
$ws.Range("E27").Value = "  +0.49%  "

# Row 28 - RenderToken
$ws.Range("D28").Value = "'9.23"
$ws.Range("E28").Value = "  -0.29%  "

# Row 29 - ImmutableX
$ws.Range("E29").Value = "  +7.35%  "

# Row 30 - Stacks
$ws.Range("D30").Value = "'3.06"
$ws.Range("E30").Value = "  +6.72%  "

# Row 31 - NEARProtocol
$ws.Range("D31").Value = "'7.01"
$ws.Range("E31").Value = "  +6.85%  "

# Row 32 - EthereumClassic
$ws.Range("D32").Value = "'28.13"
$ws.Range("E32").Value = "  +0.27%  "

# Row 33 - Mantle
$ws.Range("E33").Value = "  -1.52%  "

# Row 34 - FirstDigitalUSD
$ws.Range("E34").Value = "  +0.12%  "

# Row 35 - Filecoin
$ws.Range("D35").Value = "'6.55"
$ws.Range("E35").Value = "  -0.71%  "

# Row 36 - Bittensor
$ws.Range("D36").Value = "'513.74"
$ws.Range("E36").Value = "  +6.20%  "

# Row 37 - OKB
$ws.Range("D37").Value = "'54.82"
$ws.Range("E37").Value = "  -0.85%  "

# Row 38 - Hedera
$ws.Range("D38").Value = "'0.0896"
$ws.Range("E38").Value = "  -2.03%  "

# Row 39 - VeChain
$ws.Range("D39").Value = "'0.0423"
$ws.Range("E39").Value = "  -0.12%  "

# Row 40 - Kaspa
$ws.Range("E40").Value = "  +6.29%  "

# Row 41 - Cosmos
$ws.Range("D41").Value = "'8.88"
$ws.Range("E41").Value = "  -0.23%  "

# Row 42 - TheGraph
$ws.Range("D42").Value = "'0.304"
$ws.Range("E42").Value = "  +4.53%  "

# Row 43 - dogwifhat (-> PEPE)
$ws.Range("B43").Value = "PEPE"
$ws.Range("C43").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D43").Value = "'0.0₃0680"
$ws.Range("E43").Value = "  +8.56%  "

# Row 44 - PEPE (-> dogwifhat)
$ws.Range("B44").Value = "dogwifhat"
$ws.Range("C44").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D44").Value = "'2.86"
$ws.Range("E44").Value = "  -4.90%  "

# Row 45 - Fetch.AI
$ws.Range("D45").Value = "'2.45"
$ws.Range("E45").Value = "  -1.43%  "

# Row 46 - Maker
$ws.Range("D46").Value = "'2.859.25"
$ws.Range("E46").Value = "  -5.23%  "

# Row 47 - InjectiveProtocol
$ws.Range("D47").Value = "'28.49"
$ws.Range("E47").Value = "  -1.51%  "

# Row 48 - ThetaToken
$ws.Range("D48").Value = "'2.40"
$ws.Range("E48").Value = "  +4.36%  "

# Row 50 - Stellar
$ws.Range("E50").Value = "  +0.91%  "

# Row 51 - CoreDAO
$ws.Range("D51").Value = "'2.62"
$ws.Range("E51").Value = "  +6.53%  "

